$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the correlation values in B2:D9 to 0, per the diff.
$ws.Range("B2:D9").Value = 0
